# fix: prevent hidden columns from being labeled upon detecting changes
#
# The change-detection logic used to build this AHB comparison sheet was
# comparing the FV2304 columns (B:K) against the FV2310 columns (M:V)
# including some hidden helper columns. That made it wrongly flag a few
# rows with the "change" marker (column L, "AENDERUNG") even though the
# two sides are actually identical for those rows.
#
# This script removes the erroneous marker from column L on the affected
# rows. For rows that are also the first row of a new segment group, the
# whole row additionally gets the shaded "group header" look (matching
# rows such as 2, 9, 12, 15, 19, 23) instead of the plain row look it
# incorrectly had.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cells whose current formatting represents the desired end
# state; row 9 is an already-correct "group header" row.
$styleRegularSrc = $ws.Range("A9")   # plain "group header" cell look
$styleBoldSrc    = $ws.Range("B9")   # bold "group header" cell look (column B)
$styleChangeSrc  = $ws.Range("L9")   # cleared "change" column look

# Rows that are the first row of a new segment group and have no real
# left/right difference: re-shade the whole row as a group header row.
$fullRestyleRows = @(30, 33, 40, 47, 51, 54, 78)

# Rows that have no real left/right difference but are not the first row
# of their group: only the erroneous "change" marker needs clearing.
$changeOnlyRows = @(31, 32, 34, 35, 36, 37, 38, 39, 41, 42, 44, 45, 48, 49, 50, 52, 53, 55, 56, 64, 65, 66, 72, 73, 74, 76, 77, 79, 80)

foreach ($row in $fullRestyleRows) {
    # Column B -> bold "group header" look.
    $styleBoldSrc.Copy()
    $ws.Range("B$row").PasteSpecial(-4122)

    # Column A and C:K, M:V -> plain "group header" look.
    $styleRegularSrc.Copy()
    $ws.Range("A$row").PasteSpecial(-4122)

    $styleRegularSrc.Copy()
    $ws.Range("C$row`:K$row").PasteSpecial(-4122)

    $styleRegularSrc.Copy()
    $ws.Range("M$row`:V$row").PasteSpecial(-4122)

    # Column L -> clear the erroneous "change" marker.
    $ws.Range("L$row").Value = ""
    $styleChangeSrc.Copy()
    $ws.Range("L$row").PasteSpecial(-4122)
}

foreach ($row in $changeOnlyRows) {
    # Column L -> clear the erroneous "change" marker.
    $ws.Range("L$row").Value = ""
    $styleChangeSrc.Copy()
    $ws.Range("L$row").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
